# Auto-generated Excel COM-interop script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.494.97'
$ws.Range("E2").Value = '  +0.94%  '

$ws.Range("D3").Value = '1.826.15'
$ws.Range("E3").Value = '  +1.28%  '

$ws.Range("D4").Value = '''1.003'
$ws.Range("E4").Value = '  +0.56%  '

$ws.Range("D5").Value = '''317.67'
$ws.Range("E5").Value = '  +0.34%  '

$ws.Range("D6").Value = '''1.002'
$ws.Range("E6").Value = '  +0.44%  '

$ws.Range("D7").Value = '''0.5346'
$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").Value = '''0.3961'
$ws.Range("E8").Value = '  +5.04%  '

$ws.Range("D9").Value = '''0.07740'
$ws.Range("E9").Value = '  +3.89%  '

$ws.Range("D10").Value = '''1.118'
$ws.Range("E10").Value = '  +2.00%  '

$ws.Range("D11").Value = '''41.98'
$ws.Range("E11").Value = '  -0.14%  '

$ws.Range("D12").Value = '''6.364'
$ws.Range("E12").Value = '  +3.72%  '

$ws.Range("D13").Value = '''21.21'
$ws.Range("E13").Value = '  +3.03%  '

$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = '''7.588'
$ws.Range("E14").Value = '  +4.16%  '

$ws.Range("B15").Value = 'BinanceUSD'
$ws.Range("C15").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D15").Value = '''1.003'
$ws.Range("E15").Value = '  +0.57%  '

$ws.Range("D16").Value = '1.823.21'
$ws.Range("E16").Value = '  +1.58%  '

$ws.Range("D17").Value = '''0.00001089'
$ws.Range("E17").Value = '  +2.91%  '

$ws.Range("D18").Value = '''89.83'
$ws.Range("E18").Value = '  +0.78%  '

$ws.Range("D19").Value = '''0.06579'
$ws.Range("E19").Value = '  +1.49%  '

$ws.Range("D20").Value = '''17.77'
$ws.Range("E20").Value = '  +3.29%  '

$ws.Range("D21").Value = '''1.002'
$ws.Range("E21").Value = '  +0.46%  '

$ws.Range("D22").Value = '''6.075'
$ws.Range("E22").Value = '  +2.85%  '

$ws.Range("D23").Value = '28.510.03'
$ws.Range("E23").Value = '  +1.03%  '

$ws.Range("D24").Value = '''11.19'
$ws.Range("E24").Value = '  +0.28%  '

$ws.Range("D25").Value = '''2.256'
$ws.Range("E25").Value = '  +8.26%  '

$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").Value = '''20.74'
$ws.Range("E26").Value = '  +2.05%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '''157.29'
$ws.Range("E27").Value = '  +1.09%  '

$ws.Range("D28").Value = '''2.437'
$ws.Range("E28").Value = '  +5.01%  '

$ws.Range("D29").Value = '2.035.36'
$ws.Range("E29").Value = '  +1.62%  '

$ws.Range("D30").Value = '''125.69'
$ws.Range("E30").Value = '  +3.76%  '

$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = '''0.1124'
$ws.Range("E31").Value = '  +5.09%  '

$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = '''1.132'
$ws.Range("E32").Value = '  +1.32%  '

$ws.Range("D33").Value = '''5.732'
$ws.Range("E33").Value = '  +2.86%  '

$ws.Range("D34").Value = '''3.660'
$ws.Range("E34").Value = '  +0.43%  '

$ws.Range("D35").Value = '''0.07315'
$ws.Range("E35").Value = '  +4.59%  '

$ws.Range("D36").Value = '''0.2259'
$ws.Range("E36").Value = '  +1.29%  '

$ws.Range("D37").Value = '''0.02356'
$ws.Range("E37").Value = '  +3.00%  '

$ws.Range("D38").Value = '''8.958'
$ws.Range("E38").Value = '  +5.79%  '

$ws.Range("D39").Value = '''5.204'
$ws.Range("E39").Value = '  +3.37%  '

$ws.Range("D40").Value = '''11.41'
$ws.Range("E40").Value = '  +2.58%  '

$ws.Range("D41").Value = '''0.6284'
$ws.Range("E41").Value = '  +1.87%  '

$ws.Range("D42").Value = '''1.199'
$ws.Range("E42").Value = '  +1.43%  '

$ws.Range("D43").Value = '''1.002'
$ws.Range("E43").Value = '  +0.44%  '

$ws.Range("D44").Value = '''1.397'
$ws.Range("E44").Value = '  -2.69%  '

$ws.Range("D45").Value = '''13.55'
$ws.Range("E45").Value = '  +2.29%  '

$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").Value = '''3.720'
$ws.Range("E46").Value = '  +1.24%  '

$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '''0.5891'
$ws.Range("E47").Value = '  +1.98%  '

$ws.Range("D48").Value = '''125.45'
$ws.Range("E48").Value = '  +0.09%  '

$ws.Range("D49").Value = '''1.997'
$ws.Range("E49").Value = '  +3.80%  '

$ws.Range("D50").Value = '''1.196'
$ws.Range("E50").Value = '  +0.18%  '

$ws.Range("D51").Value = '''0.06930'
$ws.Range("E51").Value = '  +1.11%  '
